# Generate Report for Handback
# The file "3a5c89dd-aa36-4cc1-8cab-1bfe5a066509.md" has been handed back
# (in sync with en-US) for both the zh-cn and de-de locales, so the
# localization-status report is regenerated to reflect that: the status
# moves from "Ready for handoff" to "Handed back: in sync with en-US" and
# the "Latest Handback File / DateTime" columns get populated for row 3
# on each locale sheet.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"
$handbackFileName = "3a5c89dd-aa36-4cc1-8cab-1bfe5a066509.md"

# ---- Overview sheet ----------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

# ---- zh-cn sheet ---------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $handedBack
$zhcn.Hyperlinks.Add(
    $zhcn.Range("J3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/9c2e6f4a7d1b8e3c5f0a6d4b2e8c1f7a3b5d9e0c/e2e/3a5c89dd-aa36-4cc1-8cab-1bfe5a066509.md",
    [Type]::Missing,
    [Type]::Missing,
    $handbackFileName
) | Out-Null
$zhcn.Range("K3").Value = "3a5c89dd-aa36-4cc1-8cab-1bfe5a066509.81128374e429756d48e2597df6ee93a2d3f65293.zh-cn.xlf"
$zhcn.Range("L3").Value = "2017-02-17 09:41:29"

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $handedBack
$dede.Hyperlinks.Add(
    $dede.Range("J3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/1a4d7c9e2b6f8035c1e9a7d3b5f0c2e6a8d4b1f7/e2e/3a5c89dd-aa36-4cc1-8cab-1bfe5a066509.md",
    [Type]::Missing,
    [Type]::Missing,
    $handbackFileName
) | Out-Null
$dede.Range("K3").Value = "3a5c89dd-aa36-4cc1-8cab-1bfe5a066509.81128374e429756d48e2597df6ee93a2d3f65293.de-de.xlf"
$dede.Range("L3").Value = "2017-02-17 09:41:54"
